# feat: add 2022-Q4 data
#
# The previously-sole "2022-Q3" sheet is renamed to "2022-Q4" and filled
# with the new quarter's fund-holdings table; a fresh "2022-Q3" sheet is
# inserted right after it, preserving the original Q3 fund-holdings table
# (via a same-workbook sheet copy, so every value/style round-trips
# byte-for-byte). The summary sheet ("总计") gets its 2022-Q3 row's data
# replaced by 2022-Q4 numbers, and a new row is appended underneath with
# the original 2022-Q3 summary figures.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate the current "2022-Q3" sheet *before* we touch it, so the
#    duplicate keeps the original Q3 fund-holdings data/styles intact.
#    Excel places the copy immediately after the source sheet.
# ------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($null, $q3Sheet)
$newQ3Sheet = $wb.Worksheets.Item($q3Sheet.Index + 1)

# ------------------------------------------------------------------
# 2. Turn the original sheet into "2022-Q4" and overwrite its data with
#    the new quarter's fund-holdings table.
# ------------------------------------------------------------------
$q4Sheet = $q3Sheet
$q4Sheet.Name = "2022-Q4"

# Re-stamp the header row + the "序号" (A) column with the same bold /
# centered / bordered style ("s=2") used on the "总计" sheet - the new
# quarter's table was authored with that style rather than the s=1 style
# the old Q3 sheet uses.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("B1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$q4Sheet.Range("A2:A16").PasteSpecial(-4122)

# Helper: write a value as literal text (keeps leading zeros / decimal
# text exactly as given) without leaving a stray "Text" number-format
# style behind on the cell.
function Set-TextValue($range, $value) {
  $range.NumberFormat = "@"
  $range.Value = $value
  $range.Style = "Normal"
}

$q4Rows = @(
  @(0,  "013389", "华夏成长先锋一年持有混合A",       "11.63", "80.24", "8.49", "0.9874", 1),
  @(1,  "013390", "华夏成长先锋一年持有混合C",       "3.59",  "80.24", "8.49", "0.3048", 1),
  @(2,  "014410", "华夏时代领航两年持有混合A",       "2.51",  "70.90", "8.98", "0.2254", 1),
  @(3,  "000800", "华商未来主题混合",                 "4.12",  "74.31", "3.86", "0.1590", 4),
  @(4,  "010656", "华商均衡30混合",                   "3.86",  "33.52", "3.75", "0.1448", 2),
  @(5,  "011371", "华商远见价值混合型证券投资基金A", "3.24",  "81.85", "4.40", "0.1426", 3),
  @(6,  "630010", "华商价值精选混合",                 "4.37",  "87.25", "3.21", "0.1403", 7),
  @(7,  "001449", "华商双驱优选灵活配置混合",         "2.26",  "77.71", "5.27", "0.1191", 1),
  @(8,  "008961", "华商科技创新混合",                 "2.56",  "93.13", "2.82", "0.0722", 9),
  @(9,  "002289", "华商改革创新股票A",                 "1.12",  "85.85", "4.07", "0.0456", 2),
  @(10, "014411", "华夏时代领航两年持有混合C",       "0.45",  "70.90", "8.98", "0.0404", 1),
  @(11, "010403", "华商景气优选混合",                 "0.58",  "77.20", "4.96", "0.0288", 2),
  @(12, "630006", "华商产业升级混合",                 "0.86",  "88.65", "3.25", "0.0280", 7),
  @(13, "016052", "华商改革创新股票C",                 "0.48",  "85.85", "4.07", "0.0195", 2),
  @(14, "011372", "华商远见价值混合型证券投资基金C", "0.29",  "81.85", "4.40", "0.0128", 3)
)

$r = 2
foreach ($row in $q4Rows) {
  $q4Sheet.Cells.Item($r, 1).Value = $row[0]
  Set-TextValue $q4Sheet.Cells.Item($r, 2) $row[1]
  Set-TextValue $q4Sheet.Cells.Item($r, 3) $row[2]
  Set-TextValue $q4Sheet.Cells.Item($r, 4) $row[3]
  Set-TextValue $q4Sheet.Cells.Item($r, 5) $row[4]
  Set-TextValue $q4Sheet.Cells.Item($r, 6) $row[5]
  Set-TextValue $q4Sheet.Cells.Item($r, 7) $row[6]
  $q4Sheet.Cells.Item($r, 8).Value = $row[7]
  $r = $r + 1
}

# ------------------------------------------------------------------
# 3. Rename the duplicated sheet back to "2022-Q3" (data is already the
#    original Q3 fund-holdings table, untouched).
# ------------------------------------------------------------------
$newQ3Sheet.Name = "2022-Q3"

# ------------------------------------------------------------------
# 4. Update the "总计" summary sheet: row 2 becomes the 2022-Q4 figures,
#    and a new row 3 is added below with the original 2022-Q3 figures.
# ------------------------------------------------------------------

# Copy row 2's formatting down into the new row 3 first (so A3 keeps the
# bold/centered/bordered style that A2 has) before changing any values.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(3,2).Value = "2022-Q3"
$totalSheet.Cells.Item(3,3).Value = 4
$totalSheet.Cells.Item(3,4).Value = 1.25

$totalSheet.Cells.Item(2,2).Value = "2022-Q4"
$totalSheet.Cells.Item(2,3).Value = 15
$totalSheet.Cells.Item(2,4).Value = 2.47
